$d = $word.ActiveDocument

# 1) Replace the lone "@" placeholder with {client_adress}
$d.Content.Find.Execute("@", $true, $false, $false, $false, $false, $true, 1, $false, "{client_adress}", 2)

# 2) Replace the lone "*" placeholder with {oa_contact}
$d.Content.Find.Execute("*", $true, $false, $false, $false, $false, $true, 1, $false, "{oa_contact}", 2)

# 3) Split "(TVA 7.7 % incluse), ..." so "TVA 7.7 % incluse" becomes its own
#    run containing "{TVA}", leaving the "(" and "), ..." runs behind with
#    matching (unchanged) run formatting.
$rng = $d.Content
$rng.Find.Execute("TVA 7.7 % incluse", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tvaRange = $d.Range($rng.Start, $rng.End)
$tvaRange.Text = "{TVA}"
# Toggling a character property forces the engine to materialize the
# edited span as its own run (then restore the original formatting so the
# run properties end up identical to their neighbours).
$tvaRange.Bold = 1
$tvaRange.Bold = 0
